# attack speed + anim speed player
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Player")

# Row 3 (Level 1)
$ws.Range("B3").Value = 11111
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 0.1

# Row 4 (Level 2)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 3

# Row 5 (Level 3)
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 3

# Row 6 (Level 4)
$ws.Range("E6").Value = 0.5
$ws.Range("F6").Value = 3

# Update the selected cell to match the author's final cursor position
$ws.Range("E13").Select()
